$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '247.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.65%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '29.53'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '8.70%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.165'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.14%'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.59%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.577'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.96%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8569'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.57%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8681'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.95%'
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.01030'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1,608.62%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1366'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.36%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07064'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.58%'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.55%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09383'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.11%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001528'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.01%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04135'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.24%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005980'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.75%'
$ws.Range("B17").Value = 'UpBots'
$ws.Range("C17").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.007489'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '5,070.82%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.489'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.58%'
$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.099'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.97%'
$ws.Range("B20").Value = 'BTSEToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.281'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.54%'
$ws.Range("B21").Value = 'BitpandaEcosystemToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3184'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.62%'
$ws.Range("B22").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C22").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.03390'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '5.37%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.28%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.463'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.79%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.42%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.005009'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '11.86%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.70%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '22.20%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03752'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.56%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005787'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-2.54%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1072'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.38%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002427'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '2.07%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008479'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-12.76%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005250'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.59%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.01%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002276'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-9.26%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.01%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
